# Reorder the "Requisitos" bullet list in LOB1223.docx.
#
# The list lives inside a *single* paragraph (style "List Bullet") made up
# of many runs, each holding one "CODE -  Name  (Requisito fraco)" line
# followed by a manual line break (w:br). The edit only reshuffles the
# existing lines -- no lines are added, removed, or reworded -- so we
# rebuild that paragraph's text in the new order, keeping the same
# manual-line-break (vertical tab, chr(11)) separators Word uses for w:br.

$d = $word.ActiveDocument

# Final, target order of the requirement lines.
$newLines = @(
    'LOB1257 -  Sistema de Abastecimento e Tratamento de Água  (Requisito fraco)',
    'LOB1019 -  Física II  (Requisito fraco)',
    'LOB1011 -  Eletricidade Aplicada  (Requisito fraco)',
    'LOB1037 -  Àlgebra Linear  (Requisito fraco)',
    'LOB1232 -  Licenciamento Ambiental  (Requisito fraco)',
    'LOB1038 -  Física Experimental I  (Requisito fraco)',
    'LOB1052 -  Cálculo III  (Requisito fraco)',
    'LOQ4233 -  Gestão de Negócios  (Requisito fraco)',
    'LOB1042 -  Física Experimental IV  (Requisito fraco)',
    'LOB1039 -  Física Experimental III  (Requisito fraco)',
    'LOQ4097 -  Fundamentos de Química para Engenharia I (Requisito fraco)',
    'LOB1024 -  Mecânica  (Requisito fraco)',
    'LOB1053 -  Física III  (Requisito fraco)',
    'LOB1041 -  Física Experimental II  (Requisito fraco)',
    'LOB1012 -  Estatística  (Requisito fraco)',
    'LOB1003 -  Cálculo I  (Requisito fraco)',
    'LOB1004 -  Cálculo II  (Requisito fraco)',
    'LOB1045 -  Leitura e Produção de Textos Acadêmicos  (Requisito fraco)',
    'LOB1006 -  Cálculo IV  (Requisito fraco)',
    'LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito fraco)',
    'LOQ4247 -  Desenho Assistido por Computador  (Requisito fraco)',
    'LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito fraco)',
    'LOB1021 -  Física IV  (Requisito fraco)',
    'LOM3081 -  Introdução à Mecânica dos Sólidos  (Requisito fraco)',
    'LOB1036 -  Geometria Analítica  (Requisito fraco)',
    'LOB1018 -  Física I  (Requisito fraco)',
    'LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito fraco)',
    'LOQ4095 -  Química Geral Experimental  (Requisito fraco)'
)

# Locate the "Requisitos" heading paragraph, then the List Bullet paragraph
# right after it -- that is the requirement list we need to rebuild.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($text -eq "Requisitos") {
        $target = $d.Paragraphs.Item($i + 1)
        break
    }
}

# Manual line break character Word uses inside Range.Text for <w:br/>.
$vbreak = [char]11
$joined = [string]::Join($vbreak, $newLines)

$start = $target.Range.Start
$end = $target.Range.End
$r = $d.Range($start, $end)
$r.Text = $joined
